# CenterDb 추가 / CookieSoulStone 추가 / Log 성 테이블 추가
#
# This workbook (Cookie.xlsx) describes the "Packet" model for the Cookie
# table. The StarExp / AccStarExp fields are being renamed to
# SoulStone / AccSoulStone.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Packet")

# Row 4: StarExp -> SoulStone
$ws.Range("A4").Value = "SoulStone"

# Row 5: AccStarExp -> AccSoulStone
$ws.Range("A5").Value = "AccSoulStone"

# Move/restore the active selection to B8 (matches the saved sheet view).
[void]$ws.Range("B8").Select()
